# Apply the change: insert one new data row at row 155 in the "Hortaliza,
# Vega Monumental Concepción - Zapallo italiano" sheet, shifting the
# existing rows 155-196 down to 156-197, and populate the new row 155
# with its own data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 155, shifting rows 155:196 down to 156:197.
# -4121 = xlShiftDown
$ws.Rows("155:155").Insert(-4121)

# Populate the newly inserted row 155 with its data (same dimension/style
# was inherited from the row above, matching column D's date style).
$ws.Range("A155").Value = 11
$ws.Range("B155").Value = "Vega Monumental Concepción"
$ws.Range("C155").Value = "Bíobío"
$ws.Range("D155").Value = 44985
$ws.Range("E155").Value = 8
$ws.Range("F155").Value = 100112032
$ws.Range("G155").Value = "Zapallo italiano"
$ws.Range("H155").Value = "Sin especificar"
$ws.Range("I155").Value = "Primera"
$ws.Range("J155").Value = 250
$ws.Range("K155").Value = 8000
$ws.Range("L155").Value = 9000
$ws.Range("M155").Value = 8400
$ws.Range("N155").Value = "$/caja 50 unidades"
$ws.Range("O155").Value = "Región de Arica y Parinacota"
$ws.Range("P155").Value = 168
$ws.Range("Q155").Value = 50
$ws.Range("R155").Value = "Hortaliza"
